$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Before:
#   P1  2022年6月1日，星期三。
#   P2  多云，今天是六一儿童节，又是开心的一天呢。
#   P3  2022年6月2日星期四
#   P4  (pPr/rFonts hint=eastAsia) 中雨，今天是农历五月初四，明天就是端午节了。
#
# After:
#   P1  2022年6月1日，星期三。                                   (unchanged)
#   P2  多云，今天是六一儿童节，又是开心的一天呢。                (unchanged)
#   P3  2022年6月2日星期四                                        (unchanged)
#   P4  中雨，今天是农历五月初四，明天就是端午节了。              (new, plain)
#   P5  2022年6月3日星期五                                        (new, plain)
#   P6  (pPr/rFonts hint=eastAsia) 中雨，今天是农历五月初五，中国传统端午节。
#   P7  (pPr/rFonts hint=eastAsia, empty)                         (new, empty)
# ---------------------------------------------------------------------------

# 1) Insert a new plain paragraph right after paragraph 3 that carries the
#    old weather/date line verbatim ("...初四，明天就是端午节了。").
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "中雨，今天是农历五月初四，明天就是端午节了。"

# 2) Insert another new plain paragraph after it for the 6/3 date line. This
#    line mixes hinted (CJK) and un-hinted (ASCII) runs just like the other
#    date paragraphs in the document, so we inject the exact run structure
#    via InsertXML instead of a single Range.Text assignment.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$dateXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r>' +
           '<w:r><w:t>022</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>年6月</w:t></w:r>' +
           '<w:r><w:t>3</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>日星期</w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>五</w:t></w:r>' +
           '</w:p>'
$p5.Range.InsertXML($dateXml)

# 3) The original fourth paragraph (now the sixth) keeps its pPr, but its
#    wording changes to describe the 5th of the lunar month / the festival
#    itself rather than tomorrow's forecast.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Find.Execute("中雨，今天是农历五月初四，明天就是端午节了。", $true, $false, $false, `
    $false, $false, $true, 1, $false, "中雨，今天是农历五月初五，中国传统端午节。", 2)

# 4) Append a trailing empty paragraph that keeps the same pPr/rPr hint but
#    has no run content at all.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$emptyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
            '</w:p>'
$p7.Range.InsertXML($emptyXml)
